# Regenerate the s_vals rows (B2:G9) with the updated values produced by
# the filtered (save-game-excluded) re-run of the data pipeline.
# Column G is the row sum of B:E, recomputed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 3.182878228561681,  1.65323645889881,    3.082599426703578,  0.4998867070740569, 0, 8.418600821238126),
    @(3, 3.182878228561681,  1.65323645889881,    0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(4, 1.505614041169197,  1.65323645889881,    0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054),
    @(5, 0.3464964993005633, 0.004309184025731883,0.7127328510149897, 0.4998867070740569, 0, 1.563425241415342),
    @(6, 3.182878228561681,  1.65323645889881,    0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729),
    @(7, 0.3464964993005633, 0.3375848360084654,  16.98373111632243,  6.48142807727062,   0, 24.14924052890208),
    @(8, 1.505614041169197,  0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1, 2.210719231951476),
    @(9, 3.182878228561681,  1.65323645889881,    3.082599426703578,  6.48142807727062,   0, 14.40014219143469)
)

foreach ($rowVals in $values) {
    $r = $rowVals[0]
    $ws.Range("B$r").Value = $rowVals[1]
    $ws.Range("C$r").Value = $rowVals[2]
    $ws.Range("D$r").Value = $rowVals[3]
    $ws.Range("E$r").Value = $rowVals[4]
    $ws.Range("F$r").Value = $rowVals[5]
    $ws.Range("G$r").Value = $rowVals[6]
}
